$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("priorities")

# Update I1 header text (comments_for_submission -> comments_for_submission_with_fish_permit)
$ws.Range('I1').Value2 = 'comments_for_submission_with_fish_permit'

# Update H column comment text for rows with new content
$ws.Range('H14').Value2 = 'Undercut banks, large woody debris and overhanging vegetation througout.  Pools shallow. Beaver dams start 330m upstream of crossing. Minnowtrapping conducted upstream and downstream of crossing with Rainbow Trout and Sculpin captured downstream. '
$ws.Range('H19').Value2 = ' Some deep pools and boulders, udercut banks,  large wody debris and gravels throughout.   Some debris steps from 30 - 70 cms high. Passble railway culvert located downstream (16603641). New bridge upstream.'
$ws.Range('H20').Value2 = 'Abundant undercut banks, overhanging vegetation, large woody debris and gravels.  Historic beaver dam 700 m upstream.  Railway culvert (modelled crossing 16603287) is  200 m upstream and is barrier (90 m long, unembedded and 3%).'
$ws.Range('H23').Value2 = 'CN Rail crossing.  Abundant gravels, large woody debris, undercut banks, overhanging vegetation and small woody debris. Recently installed bridges downsteam and upstream.  20 cm long bull trout (suspected) observed approximately 340 m upstream of the culvert.  Minnowtrapping conducted upstream and downstream with Rainbow Trout captured downstream.'
$ws.Range('H25').Value2 = 'Stable channel with large woody debris throughout.  Railway crossing culvert (modelled ID 16603267) is located 60 m downstream of the crossing and is a barrier. Overhanging vegetation and undercut banks present for cover. Historic beaver impounded area at top of site.'
$ws.Range('H26').Value2 = 'Culvert is under Chuchinka-Colbourne FSR but CN railway crossing (PSCIS 57687) is  located 10 m upstream and also has barrier crossing. Abundant gravels throughout with deep pools suitable for overwintering.'
$ws.Range('H27').Value2 = 'Salmonids and cyprinids observed downstream throughout. Beaver activity with breached dam present.  Hunting/fishing camp located just downstream of crossing near confluence of Parsnip River. '

# Upgrade CV1 rows (22, 23) to High priority
$ws.Range('G22').Value2 = 'High'
$ws.Range('G23').Value2 = 'High'

# Apply AutoFilter on priority column G, filtering to only show High priority rows
$ws.Range('G1:G38').AutoFilter(1, @('High'))

# Update the active selection/view to match the target state
$ws.Range('B25').Select()
